$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$replacements = @(
    @{ Row = 1;  Col = 1; Old = "75÷6=12, 3"; New = "11÷2=5, 1" },
    @{ Row = 1;  Col = 2; Old = "66÷5=13, 1"; New = "99÷9=11, 0" },
    @{ Row = 1;  Col = 3; Old = "12÷2=6, 0";  New = "47÷4=11, 3" },
    @{ Row = 1;  Col = 4; Old = "94÷6=15, 4"; New = "54÷3=18, 0" },
    @{ Row = 1;  Col = 5; Old = "82÷3=27, 1"; New = "51÷3=17, 0" },

    @{ Row = 5;  Col = 1; Old = "73÷7=10, 3"; New = "15÷2=7, 1" },
    @{ Row = 5;  Col = 2; Old = "52÷8=6, 4";  New = "83÷5=16, 3" },
    @{ Row = 5;  Col = 3; Old = "87÷3=29, 0"; New = "19÷6=3, 1" },
    @{ Row = 5;  Col = 4; Old = "55÷9=6, 1";  New = "24÷3=8, 0" },
    @{ Row = 5;  Col = 5; Old = "77÷5=15, 2"; New = "13÷2=6, 1" },

    @{ Row = 9;  Col = 1; Old = "82÷5=16, 2"; New = "85÷7=12, 1" },
    @{ Row = 9;  Col = 2; Old = "49÷7=7, 0";  New = "76÷8=9, 4" },
    @{ Row = 9;  Col = 3; Old = "61÷5=12, 1"; New = "48÷2=24, 0" },
    @{ Row = 9;  Col = 4; Old = "97÷9=10, 7"; New = "28÷3=9, 1" },
    @{ Row = 9;  Col = 5; Old = "50÷4=12, 2"; New = "15÷7=2, 1" },

    @{ Row = 13; Col = 1; Old = "88÷3=29, 1"; New = "56÷4=14, 0" },
    @{ Row = 13; Col = 2; Old = "18÷4=4, 2";  New = "30÷6=5, 0" },
    @{ Row = 13; Col = 3; Old = "34÷8=4, 2";  New = "40÷3=13, 1" },
    @{ Row = 13; Col = 4; Old = "39÷5=7, 4";  New = "50÷9=5, 5" },
    @{ Row = 13; Col = 5; Old = "12÷9=1, 3";  New = "23÷8=2, 7" },

    @{ Row = 17; Col = 1; Old = "76÷8=9, 4";  New = "31÷5=6, 1" },
    @{ Row = 17; Col = 2; Old = "72÷6=12, 0"; New = "97÷8=12, 1" },
    @{ Row = 17; Col = 3; Old = "17÷7=2, 3";  New = "16÷3=5, 1" },
    @{ Row = 17; Col = 4; Old = "19÷6=3, 1";  New = "56÷9=6, 2" },
    @{ Row = 17; Col = 5; Old = "83÷6=13, 5"; New = "43÷3=14, 1" }
)

foreach ($r in $replacements) {
    $cell = $tbl.Cell($r.Row, $r.Col)
    $rng = $cell.Range
    $rng.MoveEnd(1, -1) | Out-Null
    if ($rng.Text -ne $r.Old) {
        throw "Mismatch at row $($r.Row) col $($r.Col): expected '$($r.Old)' got '$($rng.Text)'"
    }
    $rng.Text = $r.New
}
